$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-04 Wednesday" "2026-02-05 Thursday"

Replace-Text "30×47=" "47×71="
Replace-Text "32×62=" "88×75="
Replace-Text "75×39=" "96×23="
Replace-Text "44×89=" "36×71="
Replace-Text "51×62=" "30×66="
Replace-Text "11×15=" "50×65="
Replace-Text "14×99=" "58×36="
Replace-Text "25×85=" "44×73="
Replace-Text "11×50=" "47×59="
Replace-Text "17×68=" "92×53="
Replace-Text "66×22=" "68×19="
Replace-Text "53×96=" "91×26="
Replace-Text "47×41=" "77×32="
Replace-Text "14×62=" "44×89="
Replace-Text "88×56=" "19×26="
Replace-Text "56×99=" "68×75="
Replace-Text "72×16=" "17×34="
Replace-Text "49×52=" "41×94="
Replace-Text "11×54=" "85×62="
Replace-Text "34×44=" "97×27="
Replace-Text "37×63=" "86×78="
Replace-Text "34×63=" "11×82="
Replace-Text "42×26=" "30×54="
Replace-Text "49×71=" "54×96="
Replace-Text "81×95=" "35×79="
